$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 3.8
$ws.Range("K2").Value = 4.5
$ws.Range("P2").Value = 2.88
$ws.Range("R2").Value = 1.76
$ws.Range("S2").Value = 2.24
$ws.Range("T2").Value = 1.52
$ws.Range("U2").Value = 2.84
$ws.Range("X2").Value = 28
$ws.Range("AA2").Value = 75
$ws.Range("AB2").Value = 16
$ws.Range("AC2").Value = 10.5
$ws.Range("AD2").Value = 16.5
$ws.Range("AE2").Value = 36
$ws.Range("AF2").Value = 16.5
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 14.5
$ws.Range("AI2").Value = 36
$ws.Range("AK2").Value = 17.5
$ws.Range("AL2").Value = 25
$ws.Range("AO2").Value = 22
$ws.Range("G3").Value = 2.6
$ws.Range("J3").Value = 1.2
$ws.Range("O3").Value = 1.01
$ws.Range("S3").Value = 1.25
$ws.Range("T3").Value = 1.79
$ws.Range("W3").Value = 1.64
$ws.Range("N4").Value = 5.5
$ws.Range("H5").Value = 4
$ws.Range("K5").Value = 3.5
$ws.Range("Q5").Value = 2.84
$ws.Range("G6").Value = 1.49
$ws.Range("K6").Value = 5.1
$ws.Range("M6").Value = 1.07
$ws.Range("P6").Value = 1.68
$ws.Range("T6").Value = 2.46
$ws.Range("W6").Value = 3
$ws.Range("R7").Value = 1.98
$ws.Range("S7").Value = 1.55
$ws.Range("J8").Value = 4.9
$ws.Range("K8").Value = 7.6
$ws.Range("O8").Value = 1.15
$ws.Range("P8").Value = 2.76
$ws.Range("R8").Value = 1.71
$ws.Range("S8").Value = 2.04
$ws.Range("V8").Value = 1.1
$ws.Range("W8").Value = 3.45
$ws.Range("AN8").Value = 5
$ws.Range("G9").Value = 3.7
$ws.Range("I9").Value = 2.42
$ws.Range("V9").Value = 1.7
$ws.Range("W9").Value = 1.37
$ws.Range("AD12").Value = 17.5
$ws.Range("AF12").Value = 1000
$ws.Range("G13").Value = 8.4
$ws.Range("H13").Value = 1.47
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 5.2
$ws.Range("R13").Value = 1.59
$ws.Range("T13").Value = 1.75
$ws.Range("U13").Value = 2.06
$ws.Range("Y13").Value = 13
$ws.Range("Z13").Value = 11
$ws.Range("AA13").Value = 16
$ws.Range("AB13").Value = 36
$ws.Range("AD13").Value = 12.5
$ws.Range("AH13").Value = 26
$ws.Range("F15").Value = 1.24
$ws.Range("H15").Value = 11
$ws.Range("I15").Value = 14.5
$ws.Range("J15").Value = 7.2
$ws.Range("K15").Value = 8.4
$ws.Range("N15").Value = 6.4
$ws.Range("R15").Value = 1.87
$ws.Range("S15").Value = 1.93
$ws.Range("P16").Value = 2.08
$ws.Range("J17").Value = 3.9
$ws.Range("K18").Value = 3.6
$ws.Range("O18").Value = 1.35
$ws.Range("S18").Value = 3.65
$ws.Range("AF18").Value = 20
$ws.Range("AK18").Value = 34
$ws.Range("L19").Value = 1.23
$ws.Range("F20").Value = 5.7
$ws.Range("P20").Value = 2.98
$ws.Range("Q20").Value = 1.4
$ws.Range("R20").Value = 1.82
$ws.Range("S20").Value = 2
$ws.Range("U20").Value = 2.38
$ws.Range("AC20").Value = 16
$ws.Range("G22").Value = 2.48
$ws.Range("H22").Value = 3.55
$ws.Range("I22").Value = 3.8
$ws.Range("J23").Value = 3.5
$ws.Range("G24").Value = 1.22
$ws.Range("I24").Value = 16.5
$ws.Range("J24").Value = 8.199999999999999
$ws.Range("N24").Value = 11
$ws.Range("O24").Value = 1.08
$ws.Range("P24").Value = 4.2
$ws.Range("Q24").Value = 1.25
$ws.Range("R24").Value = 2.3
$ws.Range("S24").Value = 1.62
$ws.Range("T24").Value = 1.63
$ws.Range("U24").Value = 2.26
$ws.Range("X24").Value = 80
$ws.Range("Y24").Value = 110
$ws.Range("Z24").Value = 230
$ws.Range("AA24").Value = 640
$ws.Range("AB24").Value = 21
$ws.Range("AC24").Value = 23
$ws.Range("AF24").Value = 13.5
$ws.Range("AG24").Value = 13.5
$ws.Range("AJ24").Value = 15
$ws.Range("AK24").Value = 12.5
$ws.Range("AL24").Value = 27
$ws.Range("AM24").Value = 120
$ws.Range("AN24").Value = 2.74
$ws.Range("P25").Value = 2.46
$ws.Range("Q25").Value = 1.67
$ws.Range("R25").Value = 1.58
$ws.Range("S25").Value = 2.66
$ws.Range("T25").Value = 1.55
$ws.Range("X25").Value = 21
$ws.Range("AF25").Value = 20
$ws.Range("AH25").Value = 14
